# Generate Report for Handback
# Refresh handback-report timestamps for the 02da8acc... file across the
# Overview sheet and the per-locale (zh-cn / de-de) detail sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: "Latest HO Xliff Generate Date" for row 2 (02da8acc...) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-04 20:55:18"

# --- zh-cn sheet: Correspond Handoff/Handback datetimes for row 2 (02da8acc...) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-04 20:55:14"
$wsZhCn.Range("K2").Value = "2016-09-04 20:55:30"

# --- de-de sheet: Correspond Handoff/Handback datetimes for row 2 (02da8acc...) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-09-04 20:55:18"
$wsDeDe.Range("K2").Value = "2016-09-04 20:55:37"
